$d = $word.ActiveDocument

# The document contains five paragraphs of the form "{{ANALISE_Cn}}." (n = 1..5).
# The trailing "." is a stray run that must be removed. For n = 4 and n = 5 the
# placeholder text itself is additionally split across several runs ("_C", "n",
# "}}") which must be merged into a single run "_Cn}}" once the "." is dropped.
# We locate each paragraph by searching for the "{{ANALISE_C" marker so the
# script does not depend on fixed paragraph indices.

for ($n = 1; $n -le 5; $n++) {
    $marker = "{{ANALISE_C" + $n + "}}."
    $found = $false
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $pRange = $p.Range
        if ($pRange.Text.StartsWith($marker)) {
            $pStart = $pRange.Start
            $dotStart = $pStart + $marker.Length - 1
            $dotEnd = $pStart + $marker.Length

            if ($n -eq 4 -or $n -eq 5) {
                # Merge the "_C", "n", "}}" runs (and drop the ".") into one run.
                $mergeStart = $pStart + ("{{ANALISE".Length)
                $mergeRange = $d.Range($mergeStart, $dotEnd)
                $mergeRange.Text = "_C" + $n + "}}"
            } else {
                # Just delete the stray "." run.
                $dotRange = $d.Range($dotStart, $dotEnd)
                $dotRange.Delete()
            }
            $found = $true
            break
        }
    }
}
